# Update countries & provincias Spain
#
# The source "Pais" sheet lists one country per row (col A) with its
# COVID stats in columns B-H (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes). A handful of small
# territories near the bottom of the table were re-sorted/re-labelled and
# four existing countries (Estados Unidos, Japon, Argentina, Barein) got
# refreshed daily figures. The "Datos actualizados..." banner timestamp
# also moved forward 30 minutes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refreshed daily totals for existing countries -------------------
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1684956
$ws.Range("C4").Value = 18128
$ws.Range("E4").Value = 1134096
$ws.Range("G4").Value = 595
$ws.Range("H4").Value = 99278

# Japon (row 43)
$ws.Range("B43").Value = 16550
$ws.Range("C43").Value = 14
$ws.Range("D43").Value = 13413
$ws.Range("E43").Value = 2317
$ws.Range("G43").Value = 12
$ws.Range("H43").Value = 820

# Argentina (row 48)
$ws.Range("D48").Value = 3732
$ws.Range("E48").Value = 7173
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 448

# Barein (row 53)
$ws.Range("E53").Value = 4537
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 14

# --- Re-sorted / re-labelled small territories (rows 198-216) --------
# Each row keeps its position but the country name (and therefore the
# stats that travel with it) shifts to match the new ordering.

# row 198: Nueva Caledonia -> Belice
$ws.Range("A198").Value = "Belice"
$ws.Range("B198").Value = 18
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 16
$ws.Range("E198").Value = 0
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 2

# row 199: Santa Lucia -> Nueva Caledonia
$ws.Range("A199").Value = "Nueva Caledonia"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 18
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# row 200: Belice -> Santa Lucia
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

# row 207: Islas Turcas y Caicos -> Groenlandia
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("B207").Value = 12
$ws.Range("C207").Value = 1
$ws.Range("D207").Value = 11
$ws.Range("E207").Value = 1
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

# row 208: Groenlandia -> Islas Turcas y Caicos
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 10
$ws.Range("E208").Value = 1
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 1

# row 210: Seychelles -> Montserrat
$ws.Range("A210").Value = "Montserrat"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 10
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1

# row 211: Montserrat -> Seychelles
$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# row 214: Sahara Occidental -> Bonaire, San Eustaquio y Saba
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B214").Value = 6
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 6
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# row 215: Bonaire, San Eustaquio y Saba -> San Bartolome
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("B215").Value = 6
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 6
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# row 216: San Bartolome -> Sahara Occidental
$ws.Range("A216").Value = "Sahara Occidental"
$ws.Range("B216").Value = 6
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 6
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

# --- Banner timestamp --------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 01:35"
